$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "United_States"
$ws.Range("E1").Value = "Other_countries"

$ws.Range("E2").Select()
